$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 53

$ws.Cells.Item($row, 1).Value = 45957
$ws.Cells.Item($row, 2).Value = "21,8014"
$ws.Cells.Item($row, 3).Value = "15,7989"
$ws.Cells.Item($row, 4).Value = "15,5362"
$ws.Cells.Item($row, 5).Value = "15,5362"

$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
